$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.449.10"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "2.629.25"
$ws.Range("E3").Value = "  -3.89%  "
$ws.Range("E4").Value = "  -0.02%  "
$s = $ws.Range("D5").Style
$ws.Range("D5").Value = "'552.48"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -2.34%  "
$s = $ws.Range("D6").Style
$ws.Range("D6").Value = "'154.55"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -3.92%  "
$ws.Range("E8").Value = "  -1.47%  "
$s = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.104"
$ws.Range("D9").Style = $s
$ws.Range("E10").Value = "  -4.32%  "
$ws.Range("E11").Value = "  -4.13%  "
$s = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.363"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "3.094.84"
$ws.Range("E13").Value = "  -3.93%  "
$s = $ws.Range("D14").Style
$ws.Range("D14").Value = "'25.80"
$ws.Range("D14").Style = $s
$ws.Range("E14").Value = "  -4.34%  "
$ws.Range("D15").Value = "62.350.09"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").Value = "2.632.40"
$ws.Range("E17").Value = "  -3.94%  "
$s = $ws.Range("D18").Style
$ws.Range("D18").Value = "'11.66"
$ws.Range("D18").Style = $s
$ws.Range("E18").Value = "  -5.32%  "
$ws.Range("E19").Value = "  -4.17%  "
$s = $ws.Range("D20").Style
$ws.Range("D20").Value = "'340.15"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -4.47%  "
$s = $ws.Range("D21").Style
$ws.Range("D21").Value = "'6.11"
$ws.Range("D21").Style = $s
$ws.Range("E21").Value = "  -7.64%  "
$s = $ws.Range("D22").Style
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +0.13%  "
$s = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.500"
$ws.Range("D23").Style = $s
$s = $ws.Range("D24").Style
$ws.Range("D24").Value = "'62.99"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("E28").Value = "  -8.18%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$s = $ws.Range("D29").Style
$ws.Range("D29").Value = "'7.12"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$s = $ws.Range("D30").Style
$ws.Range("D30").Value = "'1.33"
$ws.Range("D30").Style = $s
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("E31").Value = "  -4.73%  "
$s = $ws.Range("D32").Style
$ws.Range("D32").Value = "'161.08"
$ws.Range("D32").Style = $s
$ws.Range("E32").Value = "  -3.82%  "
$s = $ws.Range("D34").Style
$ws.Range("D34").Value = "'4.75"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -4.29%  "
$s = $ws.Range("D36").Style
$ws.Range("D36").Value = "'19.22"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -4.21%  "
$ws.Range("E37").Value = "  -4.02%  "
$s = $ws.Range("D38").Style
$ws.Range("D38").Value = "'335.64"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  -3.67%  "
$s = $ws.Range("D39").Style
$ws.Range("D39").Value = "'6.14"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  -2.25%  "
$s = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.913"
$ws.Range("D40").Style = $s
$ws.Range("E40").Value = "  -6.81%  "
$s = $ws.Range("D41").Style
$ws.Range("D41").Value = "'3.92"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -3.84%  "
$s = $ws.Range("D42").Style
$ws.Range("D42").Value = "'37.92"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("E43").Value = "  -6.33%  "
$ws.Range("E44").Value = "  +0.02%  "
$s = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.611"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  -3.14%  "
$s = $ws.Range("D46").Style
$ws.Range("D46").Value = "'19.76"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -5.86%  "
$s = $ws.Range("D47").Style
$ws.Range("D47").Value = "'10.98"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -0.83%  "
$s = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.0548"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -6.23%  "
$s = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.0962"
$ws.Range("D49").Style = $s
$ws.Range("E49").Value = "  -3.37%  "
$s = $ws.Range("D50").Style
$ws.Range("D50").Value = "'127.37"
$ws.Range("D50").Style = $s
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("E51").Value = "  -5.36%  "
